$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new resale-number row (row 24), mirroring the layout of the
# preceding rows: columns A-D hold plain text (date/time/weekday/week),
# columns E-T hold numbers.
$row = 24

# Date/Week look numeric to Excel's auto-detection, so force text entry
# (NumberFormat "@") then strip the format back off so no stray style
# sticks to the cell - matches the unstyled text cells used elsewhere in
# the sheet.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2023-06-05"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value = "21:48:00"
$ws.Cells.Item($row, 3).Value = "Monday"

$weekCell = $ws.Cells.Item($row, 4)
$weekCell.NumberFormat = "@"
$weekCell.Value = "23"
$weekCell.ClearFormats()

$ws.Cells.Item($row, 5).Value = 120443
$ws.Cells.Item($row, 6).Value = 134122
$ws.Cells.Item($row, 7).Value = 159861
$ws.Cells.Item($row, 8).Value = 130150
$ws.Cells.Item($row, 9).Value = 175014
$ws.Cells.Item($row, 10).Value = 112632
$ws.Cells.Item($row, 11).Value = 200186
$ws.Cells.Item($row, 12).Value = 219654
$ws.Cells.Item($row, 13).Value = 172373
$ws.Cells.Item($row, 14).Value = 119375
$ws.Cells.Item($row, 15).Value = 38317
$ws.Cells.Item($row, 16).Value = 34665
$ws.Cells.Item($row, 17).Value = 50353
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36700
$ws.Cells.Item($row, 20).Value = -1
